$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.313941333333333
$ws.Range("H2").Value = 15.941824
$ws.Range("I2").Value = 0.176869630377001
$ws.Range("J2").Value = 0.176869630377001
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2936666666666667
$ws.Range("N2").Value = 0.881
$ws.Range("O2").Value = 0.009113820319201367
$ws.Range("P2").Value = 0.009113820319201367
$ws.Range("Q2").Value = 1.560527438222222
$ws.Range("R2").Value = 14.044746944
$ws.Range("S2").Value = 0.001611958031179546
$ws.Range("T2").Value = 0.001611958031179546
$ws.Range("G3").Value = 5.313941333333333
$ws.Range("H3").Value = 15.941824
$ws.Range("I3").Value = 0.176869630377001
$ws.Range("J3").Value = 0.176869630377001
$ws.Range("O3").Value = 0.870405726797791
$ws.Range("P3").Value = 0.870405726797791
$ws.Range("Q3").Value = 149.0365150377173
$ws.Range("R3").Value = 1341.328635339456
$ws.Range("S3").Value = 0.1539483391767502
$ws.Range("T3").Value = 0.1539483391767502
$ws.Range("G4").Value = 5.313941333333333
$ws.Range("H4").Value = 15.941824
$ws.Range("I4").Value = 0.176869630377001
$ws.Range("J4").Value = 0.176869630377001
$ws.Range("O4").Value = 0.1204804528830076
$ws.Range("P4").Value = 0.1204804528830076
$ws.Range("Q4").Value = 20.62944472333511
$ws.Range("R4").Value = 185.665002510016
$ws.Range("S4").Value = 0.02130933316907123
$ws.Range("T4").Value = 0.02130933316907123
$ws.Range("I5").Value = 0.5461014638447835
$ws.Range("J5").Value = 0.5461014638447835
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2936666666666667
$ws.Range("N5").Value = 0.881
$ws.Range("O5").Value = 0.009113820319201367
$ws.Range("P5").Value = 0.009113820319201367
$ws.Range("Q5").Value = 4.818273869666668
$ws.Range("R5").Value = 43.36446482700001
$ws.Range("S5").Value = 0.004977070617534198
$ws.Range("T5").Value = 0.004977070617534198
$ws.Range("I6").Value = 0.5461014638447835
$ws.Range("J6").Value = 0.5461014638447835
$ws.Range("O6").Value = 0.870405726797791
$ws.Range("P6").Value = 0.870405726797791
$ws.Range("S6").Value = 0.4753298415431564
$ws.Range("T6").Value = 0.4753298415431564
$ws.Range("I7").Value = 0.5461014638447835
$ws.Range("J7").Value = 0.5461014638447835
$ws.Range("O7").Value = 0.1204804528830076
$ws.Range("P7").Value = 0.1204804528830076
$ws.Range("S7").Value = 0.0657945516840929
$ws.Range("T7").Value = 0.0657945516840929
$ws.Range("I8").Value = 0.2770289057782155
$ws.Range("J8").Value = 0.2770289057782155
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2936666666666667
$ws.Range("N8").Value = 0.881
$ws.Range("O8").Value = 0.009113820319201367
$ws.Range("P8").Value = 0.009113820319201367
$ws.Range("Q8").Value = 2.444236513222222
$ws.Range("R8").Value = 21.998128619
$ws.Range("S8").Value = 0.002524791670487622
$ws.Range("T8").Value = 0.002524791670487622
$ws.Range("I9").Value = 0.2770289057782155
$ws.Range("J9").Value = 0.2770289057782155
$ws.Range("O9").Value = 0.870405726797791
$ws.Range("P9").Value = 0.870405726797791
$ws.Range("S9").Value = 0.2411275460778845
$ws.Range("T9").Value = 0.2411275460778845
$ws.Range("I10").Value = 0.2770289057782155
$ws.Range("J10").Value = 0.2770289057782155
$ws.Range("O10").Value = 0.1204804528830076
$ws.Range("P10").Value = 0.1204804528830076
$ws.Range("S10").Value = 0.03337656802984344
$ws.Range("T10").Value = 0.03337656802984344
